# Fruta / hortaliza, semanal
# A new week of price data (2021-10-20, serial 44489) is inserted at the top
# of the "Chirimoya" price table, pushing the existing rows down by three.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 blank rows at row 59; everything previously at 59:77 shifts to 62:80.
$ws.Rows("59:61").Insert()

# Row 59 - Especial, Provincia del Elquí
$ws.Cells.Item(59,1).Value = 8
$ws.Cells.Item(59,2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(59,3).Value = "Coquimbo"
$ws.Cells.Item(59,4).Value = 44489
$ws.Cells.Item(59,5).Value = 4
$ws.Cells.Item(59,6).Value = "Fruta"
$ws.Cells.Item(59,7).Value = 100107
$ws.Cells.Item(59,8).Value = "Otros"
$ws.Cells.Item(59,9).Value = 100107002
$ws.Cells.Item(59,10).Value = "Chirimoya"
$ws.Cells.Item(59,11).Value = "Cultivar IV Región"
$ws.Cells.Item(59,12).Value = "Especial"
$ws.Cells.Item(59,13).Value = 300
$ws.Cells.Item(59,14).Value = 2200
$ws.Cells.Item(59,15).Value = 2300
$ws.Cells.Item(59,16).Value = 2250
$ws.Cells.Item(59,17).Value = "$/kilo (en caja de 15 kilos)"
$ws.Cells.Item(59,18).Value = "Provincia del Elquí"
$ws.Cells.Item(59,19).Value = 2250
$ws.Cells.Item(59,20).Value = 1

# Row 60 - Primera, Provincia del Elquí
$ws.Cells.Item(60,1).Value = 8
$ws.Cells.Item(60,2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(60,3).Value = "Coquimbo"
$ws.Cells.Item(60,4).Value = 44489
$ws.Cells.Item(60,5).Value = 4
$ws.Cells.Item(60,6).Value = "Fruta"
$ws.Cells.Item(60,7).Value = 100107
$ws.Cells.Item(60,8).Value = "Otros"
$ws.Cells.Item(60,9).Value = 100107002
$ws.Cells.Item(60,10).Value = "Chirimoya"
$ws.Cells.Item(60,11).Value = "Cultivar IV Región"
$ws.Cells.Item(60,12).Value = "Primera"
$ws.Cells.Item(60,13).Value = 400
$ws.Cells.Item(60,14).Value = 1900
$ws.Cells.Item(60,15).Value = 2000
$ws.Cells.Item(60,16).Value = 1950
$ws.Cells.Item(60,17).Value = "$/kilo (en caja de 15 kilos)"
$ws.Cells.Item(60,18).Value = "Provincia del Elquí"
$ws.Cells.Item(60,19).Value = 1950
$ws.Cells.Item(60,20).Value = 1

# Row 61 - Segunda, Provincia del Elquí
$ws.Cells.Item(61,1).Value = 8
$ws.Cells.Item(61,2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(61,3).Value = "Coquimbo"
$ws.Cells.Item(61,4).Value = 44489
$ws.Cells.Item(61,5).Value = 4
$ws.Cells.Item(61,6).Value = "Fruta"
$ws.Cells.Item(61,7).Value = 100107
$ws.Cells.Item(61,8).Value = "Otros"
$ws.Cells.Item(61,9).Value = 100107002
$ws.Cells.Item(61,10).Value = "Chirimoya"
$ws.Cells.Item(61,11).Value = "Cultivar IV Región"
$ws.Cells.Item(61,12).Value = "Segunda"
$ws.Cells.Item(61,13).Value = 340
$ws.Cells.Item(61,14).Value = 1400
$ws.Cells.Item(61,15).Value = 1500
$ws.Cells.Item(61,16).Value = 1450
$ws.Cells.Item(61,17).Value = "$/kilo (en caja de 15 kilos)"
$ws.Cells.Item(61,18).Value = "Provincia del Elquí"
$ws.Cells.Item(61,19).Value = 1450
$ws.Cells.Item(61,20).Value = 1
